# Update "Estado de Cuenta" worker data: replace EDGAR JOSE MEJIA GOMEZ rows
# with JUAN DANIEL YEPES CORTINAS / HENRY RODRIGUEZ CARABALLO, refresh
# salaries, add a new "2508" period for the two remaining workers, and
# update the worker count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker count (Cant. Trabajadores) 3 -> 2
$ws.Range("C13").Value = 2

# Row 16: was CC / 1043961862 / EDGAR JOSE MEJIA GOMEZ / 2507 / 56940 / 1160000
#         now CC / 1007170466 / JUAN DANIEL YEPES CORTINAS / 2507 / 56940 / 1423500
$ws.Range("C16").Value = "1007170466"
$ws.Range("D16").Value = "JUAN DANIEL YEPES CORTINAS"
$ws.Range("E16").Value = "2507"
$ws.Range("G16").Value = 1423500

# Row 17: was CC / 1043961862 / EDGAR JOSE MEJIA GOMEZ / 2505 / 56940 / 1160000
#         now CC / 1137219777 / HENRY RODRIGUEZ CARABALLO / 2507 / 56940 / 1423500
$ws.Range("C17").Value = "1137219777"
$ws.Range("D17").Value = "HENRY RODRIGUEZ CARABALLO"
$ws.Range("E17").Value = "2507"
$ws.Range("G17").Value = 1423500

# Row 18: was CC / 1007170466 / JUAN DANIEL YEPES CORTINAS / 2507 / 56940 / 1423500
#         now CC / 1007170466 / JUAN DANIEL YEPES CORTINAS / 2508 / 56940 / 1423500
$ws.Range("E18").Value = "2508"

# Row 19: was CC / 1137219777 / HENRY RODRIGUEZ CARABALLO / 2507 / 56940 / 1423500
#         now CC / 1137219777 / HENRY RODRIGUEZ CARABALLO / 2508 / 56940 / 1423500
$ws.Range("E19").Value = "2508"
